$wb = $excel.ActiveWorkbook

# --- Add the new "Work" worksheet, placed after the last existing sheet ---
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "Work"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Move($null, $lastSheet)

# NOTE: after .Move(), the original $newSheet reference can reseat to whatever
# sheet now occupies its old position, so re-fetch "Work" by name to be safe.
$ws = $wb.Worksheets.Item("Work")

# --- Populate header row + data rows ---
# Order of first-use below intentionally mirrors the target shared-string order.
$ws.Range("A2").Value = "Transit"
$ws.Range("B1").Value = "Start X"
$ws.Range("C1").Value = "Start Y"
$ws.Range("D1").Value = "End X"
$ws.Range("E1").Value = "End Y"
$ws.Range("A3").Value = "Plot"
$ws.Range("F1").Value = "Adj Start X"
$ws.Range("G1").Value = "Adj Start Y"
$ws.Range("H1").Value = "Adj End X"
$ws.Range("I1").Value = "Adj End Y"
$ws.Range("J1").Value = "Diff Start X"
$ws.Range("K1").Value = "Diff Start Y"
$ws.Range("L1").Value = "Diff End X"
$ws.Range("M1").Value = "Diff End Y"
$ws.Range("N1").Value = "Next X"
$ws.Range("O1").Value = "Next Y"
$ws.Range("P1").Value = "Next End X"
$ws.Range("Q1").Value = "Next End Y"
$ws.Range("A1").Value = "Type"
$ws.Range("A4").Value = "Plot"

# Row 2 values
$ws.Range("B2").Value = 1301.75
$ws.Range("C2").Value = 609.6
$ws.Range("D2").Value = 1301.75
$ws.Range("E2").Value = 158.48500000000001
$ws.Range("F2").Value = 1301.75
$ws.Range("G2").Value = 609.6
$ws.Range("H2").Value = 1301.75
$ws.Range("I2").Value = 160.07249999999999
$ws.Range("J2").Formula = "=F2-B2"
$ws.Range("K2").Formula = "=G2-C2"
$ws.Range("L2").Formula = "=H2-D2"
$ws.Range("M2").Formula = "=I2-E2"
$ws.Range("N2").Value = 1301.75
$ws.Range("O2").Value = 160.07249999999999
$ws.Range("P2").Value = 1428.4849999999999
$ws.Range("Q2").Value = 160.07249999999999

# Row 3 values
$ws.Range("B3").Value = 1301.75
$ws.Range("C3").Value = 158.48500000000001
$ws.Range("D3").Value = 1428.4849999999999
$ws.Range("E3").Value = 158.48500000000001
$ws.Range("F3").Value = 1301.75
$ws.Range("G3").Value = 160.07249999999999
$ws.Range("H3").Value = 1426.8974599999999
$ws.Range("I3").Value = 160.07249999999999
$ws.Range("J3").Formula = "=F3-B3"
$ws.Range("K3").Formula = "=G3-C3"
$ws.Range("L3").Formula = "=H3-D3"
$ws.Range("M3").Formula = "=I3-E3"

# Row 4 values
$ws.Range("B4").Value = 1428.4849999999999
$ws.Range("C4").Value = 158.48500000000001
$ws.Range("D4").Value = 1428.4849999999999
$ws.Range("E4").Value = 609.6
$ws.Range("F4").Value = 1426.8974599999999
$ws.Range("G4").Value = 160.07249999999999
$ws.Range("H4").Value = 1426.8974599999999
$ws.Range("I4").Value = 609.6
$ws.Range("J4").Formula = "=F4-B4"
$ws.Range("K4").Formula = "=G4-C4"
$ws.Range("L4").Formula = "=H4-D4"
$ws.Range("M4").Formula = "=I4-E4"

# --- Freeze header row and set the selection on the new sheet ---
$ws.Range("A2").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
$ws.Range("H4").Select()

# --- "Work" becomes the tab-selected / active sheet ---
$ws.Select()
